$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("N10").Value = 0

$ws.Range("H17").Value = 1197.4
$ws.Range("J17").Value = 1197.4
$ws.Range("L17").Value = 3592.2
$ws.Range("N17").Value = -3928.2

$ws.Range("H55").Value = 391.77777
$ws.Range("I55").Value = 396.54544
$ws.Range("J55").Value = 384.2857
$ws.Range("K55").Value = 396.54544
$ws.Range("L55").Value = 384.2857
$ws.Range("M55").Value = -182.54544
$ws.Range("N55").Value = -812.2857

$ws.Range("H86").Value = 2500.6
$ws.Range("I86").Value = 2001.5
$ws.Range("K86").Value = 2001.5
$ws.Range("M86").Value = -878.5

$ws.Range("H89").Value = 2500.6
$ws.Range("I89").Value = 2001.5
$ws.Range("K89").Value = 10007.5
$ws.Range("M89").Value = -4391.5

$ws.Range("H98").Value = 361.61905
$ws.Range("I98").Value = 370.7
$ws.Range("K98").Value = 370.7
$ws.Range("M98").Value = 1127.3

$ws.Range("H103").Value = 3277.7778
$ws.Range("I103").Value = 3214.2856
$ws.Range("J103").Value = 3500
$ws.Range("K103").Value = 9642.856800000001
$ws.Range("L103").Value = 10500
$ws.Range("M103").Value = -9056.856800000001
$ws.Range("N103").Value = -11672

$ws.Range("H113").Value = 3795.8333
$ws.Range("J113").Value = 4799.6665
$ws.Range("L113").Value = 4799.6665
$ws.Range("N113").Value = -11307.6665

$ws.Range("H122").Value = 361.61905
$ws.Range("I122").Value = 370.7
$ws.Range("K122").Value = 1112.1
$ws.Range("M122").Value = 1337.9

$ws.Range("H132").Value = 9221.75
$ws.Range("I132").Value = 9436.666999999999
$ws.Range("J132").Value = 8577
$ws.Range("K132").Value = 28310.001
$ws.Range("L132").Value = 25731
$ws.Range("M132").Value = -25780.001
$ws.Range("N132").Value = -30791

$ws.Range("H137").Value = 1517.4445
$ws.Range("I137").Value = 1179.5
$ws.Range("K137").Value = 3538.5
$ws.Range("M137").Value = -988.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 3867.3333
$ws.Range("I8").Value = 3818.3333
$ws.Range("K8").Value = 3818.3333
$ws.Range("M8").Value = -3674.3333

$ws.Range("H45").Value = 3613.3333
$ws.Range("I45").Value = 1995
$ws.Range("K45").Value = 1995
$ws.Range("M45").Value = -1618

$ws.Range("H97").Value = 454.53333
$ws.Range("J97").Value = 72
$ws.Range("L97").Value = 72
$ws.Range("N97").Value = -1064

$ws.Range("H110").Value = 1443.5
$ws.Range("J110").Value = 2002
$ws.Range("L110").Value = 2002
$ws.Range("N110").Value = -6092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2613.2
$ws.Range("I20").Value = 2234
$ws.Range("J20").Value = 2992.4
$ws.Range("K20").Value = 2234
$ws.Range("L20").Value = 2992.4
$ws.Range("M20").Value = -1987
$ws.Range("N20").Value = -3486.4

$ws.Range("H29").Value = 693.6667
$ws.Range("I29").Value = 472.4
$ws.Range("J29").Value = 1800
$ws.Range("K29").Value = 472.4
$ws.Range("L29").Value = 1800
$ws.Range("M29").Value = -183.4
$ws.Range("N29").Value = -2378

$ws.Range("H86").Value = 5908.48
$ws.Range("J86").Value = 6682.154
$ws.Range("L86").Value = 6682.154
$ws.Range("N86").Value = -8928.154

$ws.Range("H89").Value = 5908.48
$ws.Range("J89").Value = 6682.154
$ws.Range("L89").Value = 33410.77
$ws.Range("N89").Value = -44642.77

$ws.Range("H105").Value = 1661.875
$ws.Range("I105").Value = 1395.8334
$ws.Range("K105").Value = 1395.8334
$ws.Range("M105").Value = 351.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3822.5334
$ws.Range("I58").Value = 3377.5557
$ws.Range("J58").Value = 4490
$ws.Range("K58").Value = 3377.5557
$ws.Range("L58").Value = 4490
$ws.Range("M58").Value = -3174.5557
$ws.Range("N58").Value = -4896

$ws.Range("H107").Value = 302.75
$ws.Range("I107").Value = 187
$ws.Range("J107").Value = 1113
$ws.Range("K107").Value = 187
$ws.Range("L107").Value = 1113
$ws.Range("M107").Value = 1733
$ws.Range("N107").Value = -4953

$ws.Range("H134").Value = 1769.2354
$ws.Range("I134").Value = 1817.3125
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 5451.9375
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -2916.9375
$ws.Range("N134").Value = -8070

$ws.Range("H136").Value = 3822.5334
$ws.Range("I136").Value = 3377.5557
$ws.Range("J136").Value = 4490
$ws.Range("K136").Value = 10132.6671
$ws.Range("L136").Value = 13470
$ws.Range("M136").Value = -7582.667099999999
$ws.Range("N136").Value = -18570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 100
$ws.Range("I33").Value = 100
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 600
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -317

$ws.Range("H39").Value = 3121
$ws.Range("J39").Value = 3631.182
$ws.Range("L39").Value = 10893.546
$ws.Range("N39").Value = -11481.546

$ws.Range("H40").Value = 124.4
$ws.Range("I40").Value = 22.7
$ws.Range("K40").Value = 90.8
$ws.Range("M40").Value = -21.8

$ws.Range("H46").Value = 1233
$ws.Range("I46").Value = 200
$ws.Range("K46").Value = 600
$ws.Range("M46").Value = -509

$ws.Range("H57").Value = 1349.7084
$ws.Range("J57").Value = 1370.5883
$ws.Range("L57").Value = 4111.7649
$ws.Range("N57").Value = -5229.7649

$ws.Range("H58").Value = 2496.25
$ws.Range("J58").Value = 2995
$ws.Range("L58").Value = 8985
$ws.Range("N58").Value = -9241

$ws.Range("H113").Value = 921.94116
$ws.Range("J113").Value = 988.5714
$ws.Range("L113").Value = 2965.7142
$ws.Range("N113").Value = -7305.7142

$ws.Range("H122").Value = 821.25
$ws.Range("I122").Value = 673.3333
$ws.Range("J122").Value = 910
$ws.Range("K122").Value = 6059.9997
$ws.Range("L122").Value = 8190
$ws.Range("M122").Value = -3609.9997
$ws.Range("N122").Value = -13090

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 373.66666
$ws.Range("I9").Value = 373.66666
$ws.Range("K9").Value = 373.66666
$ws.Range("M9").Value = -203.66666

$ws.Range("H10").Value = 10666.667
$ws.Range("I10").Value = 13500
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 13500
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = -13331
$ws.Range("N10").Value = -5338

$ws.Range("H14").Value = 27668.334
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 41002.5
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 41002.5
$ws.Range("M14").Value = -832
$ws.Range("N14").Value = -41338.5

$ws.Range("H59").Value = 8000
$ws.Range("J59").Value = 8000
$ws.Range("L59").Value = 8000
$ws.Range("N59").Value = -9166

$ws.Range("H70").Value = 2709.5
$ws.Range("J70").Value = 2419
$ws.Range("L70").Value = 2419
$ws.Range("N70").Value = -2959

$ws.Range("H73").Value = 2709.5
$ws.Range("J73").Value = 2419
$ws.Range("L73").Value = 2419
$ws.Range("N73").Value = -4291

$ws.Range("H102").Value = 1951.4117
$ws.Range("I102").Value = 1545
$ws.Range("K102").Value = 1545
$ws.Range("M102").Value = 77

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5205.5
$ws.Range("I40").Value = 5205.5
$ws.Range("K40").Value = 5205.5
$ws.Range("M40").Value = -5069.5

$ws.Range("H136").Value = 3099.5
$ws.Range("I136").Value = 3099.5
$ws.Range("K136").Value = 9298.5
$ws.Range("M136").Value = -6748.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9500
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 10375
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 10375
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -11623

$ws.Range("H65").Value = 9500
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 10375
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 51875
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -58115

$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992

$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("M96").Value = 1000
$ws.Range("N96").Value = -3746

$ws.Range("H113").Value = 694
$ws.Range("I113").Value = 739.75
$ws.Range("J113").Value = 633
$ws.Range("K113").Value = 2219.25
$ws.Range("L113").Value = 1899
$ws.Range("M113").Value = -49.25
$ws.Range("N113").Value = -6239

$ws.Range("H132").Value = 2484.7144
$ws.Range("I132").Value = 2398.8333
$ws.Range("K132").Value = 7196.499899999999
$ws.Range("M132").Value = -4666.499899999999
